# Update calibration estimates with new std-moving-average based values.
$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.05077876085555474
$wsLinear.Range("B3").Value = 0.1921900119814362
$wsLinear.Range("B4").Value = 49.59597251600847
$wsLinear.Range("B5").Value = "[1.0, 0.209838494272853, 0.06607346591695039, 0.08498369905955168, 0.07492847809770685, 0.06190681670077654, 0.24426055155572682, 0.3435636722547642, 0.2080002767046351, 0.06395890061000975, 0.017844567200195734, 0.06326386710418491, 0.068613241361313, 0.20003420259515575, 0.33146322200841605, 0.21340674208965174, 0.054165879043338507, 0.04483562597631995, 0.04947750553422725, 0.040875471039267204]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.06887512307874852
$wsNonLinear.Range("B5").Value = 0.1046406578841975
$wsNonLinear.Range("B6").Value = 51.57149564516025
$wsNonLinear.Range("B7").Value = 0.3653526503391847
$wsNonLinear.Range("B8").Value = -0.3750962047132627
$wsNonLinear.Range("B9").Value = 47.78371793903636
$wsNonLinear.Range("B10").Value = "[1.0, 0.20826516817164292, 0.06712775852713737, 0.08663779142785075, 0.07666506926832543, 0.06305255376050972, 0.24361991692219173, 0.340960039104656, 0.20784881194968924, 0.0651468894481849, 0.019431448519720063, 0.06488107482122367, 0.07006047755816869, 0.1994831505400762, 0.32881872335377466, 0.21215501907350934, 0.05550020248074177, 0.04690243670826741, 0.05139791664976233, 0.04220361940047164]"
